$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# 1) Remove the "Diagrama de Negócio" requirement row (row 19) - shifts rows 20-21 up
$ws.Rows.Item(19).Delete()

# 2) The two rows that shifted up (now rows 19 and 20) inherited the special "last row"
#    formatting that used to live on the old rows 20/21. Re-normalize columns A and C
#    on those rows to match the rest of the table (copy format from row 18).
$ws.Range("A18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Mark remaining open items as completed
$ws.Range("H8").Value2 = "Finalizada"
$ws.Range("H16").Value2 = "Finalizada"
$ws.Range("H18").Value2 = "Finalizada"
$ws.Range("H19").Value2 = "Finalizada"
$ws.Range("H20").Value2 = "Finalizada"

# 4) Burndown sheet: the Sprint 3 "Pontos Realizados" cell loses its stray calculated-column
#    number format override, matching the other (plain) cells in the column. The SUMIF formula
#    itself recalculates automatically once the backlog data changed above.
$wsB = $wb.Worksheets.Item("Burndown")
$wsB.Range("D4").Copy()
$wsB.Range("D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5) Leave the cursor where the author left it when saving
$ws.Activate()
$ws.Range("H18").Select()
